# Adressbook.xlsx — fix Diane Fournier's phone number typo and update the
# current selection on Sheet1 (matches the author's manual edit: cell E4
# retyped from "418-455-1155" to "418-455-1115", cursor left on E5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet          # Sheet1 is the tabSelected/active sheet
$ws.Activate()

$ws.Range("E4").Value = "418-455-1115"

$ws.Range("E5").Select()
